$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.290.67"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.653.21"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.13"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.18"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.59"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.117.90"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  +6.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.230.33"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000146"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.656.37"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.69"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.76"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.48"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.90"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.528"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.03"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.47"
$ws.Range("E27").Value = "  +5.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("E28").Value = "  +8.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0821"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  +6.20%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.48"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.07"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").Value = "  +10.65%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  +11.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.34"
$ws.Range("E36").Value = "  +8.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.71"
$ws.Range("E37").Value = "  +5.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "343.54"
$ws.Range("E38").Value = "  +10.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.14"
$ws.Range("E39").Value = "  +5.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.911"
$ws.Range("E40").Value = "  +8.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.45"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.33"
$ws.Range("E42").Value = "  +6.65%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0580"
$ws.Range("E43").Value = "  +5.44%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "137.18"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.10"
$ws.Range("E45").Value = "  +4.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.628"
$ws.Range("E46").Value = "  +3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.46"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0251"
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1000"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.995"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.094.89"
$ws.Range("E51").Value = "  +3.00%  "

# Reset number format on price cells back to General/Normal so only the
# cell values change (matches the source data which has no explicit style).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
